$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Set the header value for BO1 ("19-aug")
$ws.Range("BO1").Value = "19-aug"

# Set numeric values for BO2:BO25
$ws.Range("BO2").Value = 93.04000000000001
$ws.Range("BO3").Value = 76.5
$ws.Range("BO4").Value = 71.81
$ws.Range("BO5").Value = 66.45999999999999
$ws.Range("BO6").Value = 65.51000000000001
$ws.Range("BO7").Value = 75.06
$ws.Range("BO8").Value = 85.15000000000001
$ws.Range("BO9").Value = 95.02
$ws.Range("BO10").Value = 98.05
$ws.Range("BO11").Value = 88.56999999999999
$ws.Range("BO12").Value = 69.64
$ws.Range("BO13").Value = 43.23
$ws.Range("BO14").Value = 30
$ws.Range("BO15").Value = 29.14
$ws.Range("BO16").Value = 30.91
$ws.Range("BO17").Value = 40.54
$ws.Range("BO18").Value = 58.18
$ws.Range("BO19").Value = 63.33
$ws.Range("BO20").Value = 63
$ws.Range("BO21").Value = 88.97
$ws.Range("BO22").Value = 97.73999999999999
$ws.Range("BO23").Value = 98.04000000000001
$ws.Range("BO24").Value = 98
$ws.Range("BO25").Value = 83.76000000000001

# Copy formatting (bold, centered, border) from BN1 header cell to BO1
$ws.Range("BN1").Copy()
$ws.Range("BO1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

